$wb = $excel.ActiveWorkbook

# Add the new "債務" (debt) worksheet.
$ws = $wb.Worksheets.Add()
$ws.Name = "債務"

# Header row (row 1)
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data row (row 2)
$ws.Range("A2").Value = 118
$ws.Range("B2").Value = "房屋貸款"
$ws.Range("C2").Value = "吳肓仁"
$ws.Range("D2").Value = "玉山銀行臺北市中山區民生東路"
$ws.Range("E2").Value = 14036373
$ws.Range("F2").Value = "97年05月07日"
$ws.Range("G2").Value = "購屋"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-03-06"
$ws.Range("K2").Value = "吳育仁"
$ws.Range("L2").Value = 1734
$ws.Range("M2").Value = "tmp476d1"
$ws.Range("N2").Value = 118

# Move the new sheet to the end (after "保險").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Move($null, $lastSheet)
